$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G53").Value = 43020101
